function Set-TextValue {
    param($ws, $cellref, $val)
    $c = $ws.Range($cellref)
    $saved = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $saved
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.757.01'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '1.812.14'
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  +0.31%  '
Set-TextValue $ws "D5" '230.66'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  +0.28%  '
Set-TextValue $ws "D8" '39.66'
$ws.Range("E8").Value = '  -9.63%  '
$ws.Range("E9").Value = '  +5.22%  '
$ws.Range("E10").Value = '  -2.38%  '
Set-TextValue $ws "D11" '0.0996'
$ws.Range("E11").Value = '  -1.21%  '
$ws.Range("D12").Value = '2.074.29'
$ws.Range("E12").Value = '  -1.35%  '
Set-TextValue $ws "D13" '11.24'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.825.44'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws "D15" '0.667'
$ws.Range("E15").Value = '  -0.58%  '
$ws.Range("E16").Value = '  -1.42%  '
$ws.Range("D17").Value = '34.747.17'
Set-TextValue $ws "D18" '69.50'
$ws.Range("E18").Value = '  -0.65%  '
$ws.Range("E19").Value = '  -1.64%  '
Set-TextValue $ws "D20" '240.47'
$ws.Range("E20").Value = '  -1.60%  '
Set-TextValue $ws "D21" '11.92'
$ws.Range("E21").Value = '  -0.80%  '
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("E23").Value = '  +0.25%  '
$ws.Range("E24").Value = '  +2.18%  '
Set-TextValue $ws "D25" '171.69'
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("E26").Value = '  -1.30%  '
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("E29").Value = '  +1.45%  '
$ws.Range("E30").Value = '  +0.33%  '
$ws.Range("E31").Value = '  +4.04%  '
Set-TextValue $ws "D32" '0.0548'
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("E33").Value = '  -2.60%  '
Set-TextValue $ws "D34" '1.26'
$ws.Range("E34").Value = '  +16.11%  '
$ws.Range("E35").Value = '  -3.09%  '
$ws.Range("E36").Value = '  +2.89%  '
Set-TextValue $ws "D37" '92.02'
$ws.Range("E37").Value = '  -4.53%  '
$ws.Range("E38").Value = '  +4.26%  '
$ws.Range("D39").Value = '1.335.79'
$ws.Range("E39").Value = '  -0.67%  '
Set-TextValue $ws "D40" '0.0192'
$ws.Range("E40").Value = '  -0.83%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws "D41" '2.47'
$ws.Range("E41").Value = '  +0.37%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws "D42" '0.967'
$ws.Range("E42").Value = '  -3.95%  '
Set-TextValue $ws "D43" '14.34'
$ws.Range("E43").Value = '  -7.25%  '
$ws.Range("E44").Value = '  -8.41%  '
$ws.Range("E45").Value = '  -4.28%  '
Set-TextValue $ws "D46" '6.23'
$ws.Range("E46").Value = '  +0.08%  '
Set-TextValue $ws "D47" '0.0514'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").Value = '2.000.09'
$ws.Range("E48").Value = '  -0.28%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  +6.77%  '
Set-TextValue $ws "D51" '98.11'
$ws.Range("E51").Value = '  -4.42%  '
